# Update FlashScore odds grid on Sheet1 to match the 2024-09-24 re-scrape.
# Each line sets one odds cell; addressed by Cells.Item(row, col) since the
# sheet is a flat data grid (row 1 = headers, data rows = matches).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 17).Value = 1.67  # Q2: 1.73 -> 1.67
$ws.Cells.Item(2, 18).Value = 2.2  # R2: 2.1 -> 2.2

# Row 5
$ws.Cells.Item(5, 11).Value = 2.88  # K5: 3 -> 2.88
$ws.Cells.Item(5, 12).Value = 12  # L5: 11 -> 12
$ws.Cells.Item(5, 17).Value = 1.5  # Q5: 1.44 -> 1.5
$ws.Cells.Item(5, 18).Value = 2.63  # R5: 2.75 -> 2.63
$ws.Cells.Item(5, 19).Value = 1.25  # S5: 1.22 -> 1.25
$ws.Cells.Item(5, 20).Value = 3.75  # T5: 4 -> 3.75
$ws.Cells.Item(5, 23).Value = 8.5  # W5: 9 -> 8.5
$ws.Cells.Item(5, 24).Value = 6  # X5: 6.5 -> 6
$ws.Cells.Item(5, 25).Value = 11  # Y5: 10 -> 11
$ws.Cells.Item(5, 26).Value = 6.5  # Z5: 7 -> 6.5
$ws.Cells.Item(5, 28).Value = 34  # AB5: 29 -> 34
$ws.Cells.Item(5, 29).Value = 17  # AC5: 19 -> 17
$ws.Cells.Item(5, 39).Value = 81  # AM5: 67 -> 81
$ws.Cells.Item(5, 46).Value = 3.75  # AT5: 4 -> 3.75
$ws.Cells.Item(5, 53).Value = 351  # BA5: 301 -> 351
$ws.Cells.Item(5, 54).Value = 301  # BB5: 251 -> 301
$ws.Cells.Item(5, 55).Value = 451  # BC5: 401 -> 451

# Row 8
$ws.Cells.Item(8, 13).Value = 1.03  # M8: 1.04 -> 1.03
$ws.Cells.Item(8, 14).Value = 15  # N8: 13 -> 15
$ws.Cells.Item(8, 17).Value = 1.7  # Q8: 1.73 -> 1.7
$ws.Cells.Item(8, 18).Value = 2.1  # R8: 2.08 -> 2.1

# Row 9
$ws.Cells.Item(9, 14).Value = 10  # N9: 9.5 -> 10

# Row 10
$ws.Cells.Item(10, 7).Value = 2.05  # G10: 2.1 -> 2.05
$ws.Cells.Item(10, 9).Value = 4.33  # I10: 4.1 -> 4.33
$ws.Cells.Item(10, 21).Value = 2.63  # U10: 2.5 -> 2.63
$ws.Cells.Item(10, 22).Value = 1.44  # V10: 1.5 -> 1.44
$ws.Cells.Item(10, 23).Value = 4.75  # W10: 5 -> 4.75
$ws.Cells.Item(10, 24).Value = 7.5  # X10: 8 -> 7.5
$ws.Cells.Item(10, 26).Value = 17  # Z10: 19 -> 17
$ws.Cells.Item(10, 34).Value = 8  # AH10: 7.5 -> 8
$ws.Cells.Item(10, 35).Value = 21  # AI10: 19 -> 21
$ws.Cells.Item(10, 37).Value = 51  # AK10: 41 -> 51
$ws.Cells.Item(10, 38).Value = 51  # AL10: 41 -> 51
$ws.Cells.Item(10, 39).Value = 67  # AM10: 51 -> 67
$ws.Cells.Item(10, 43).Value = 41  # AQ10: 51 -> 41
$ws.Cells.Item(10, 47).Value = 11  # AU10: 10 -> 11
$ws.Cells.Item(10, 50).Value = 6  # AX10: 5.5 -> 6
$ws.Cells.Item(10, 51).Value = 29  # AY10: 26 -> 29
$ws.Cells.Item(10, 54).Value = 201  # BB10: 151 -> 201

# Row 11
$ws.Cells.Item(11, 9).Value = 3.5  # I11: 3.6 -> 3.5
$ws.Cells.Item(11, 13).Value = 1.11  # M11: 1.13 -> 1.11
$ws.Cells.Item(11, 14).Value = 6.5  # N11: 6 -> 6.5
$ws.Cells.Item(11, 35).Value = 15  # AI11: 17 -> 15
$ws.Cells.Item(11, 36).Value = 13  # AJ11: 15 -> 13

# Row 12
$ws.Cells.Item(12, 7).Value = 3.8  # G12: 4 -> 3.8
$ws.Cells.Item(12, 8).Value = 3.1  # H12: 3.2 -> 3.1
$ws.Cells.Item(12, 9).Value = 2.05  # I12: 2 -> 2.05
$ws.Cells.Item(12, 12).Value = 2.88  # L12: 2.75 -> 2.88
$ws.Cells.Item(12, 13).Value = 1.1  # M12: 1.08 -> 1.1
$ws.Cells.Item(12, 14).Value = 7  # N12: 8 -> 7
$ws.Cells.Item(12, 15).Value = 1.44  # O12: 1.4 -> 1.44
$ws.Cells.Item(12, 16).Value = 2.63  # P12: 2.75 -> 2.63
$ws.Cells.Item(12, 17).Value = 2.35  # Q12: 2.3 -> 2.35
$ws.Cells.Item(12, 18).Value = 1.57  # R12: 1.6 -> 1.57
$ws.Cells.Item(12, 23).Value = 9  # W12: 9.5 -> 9
$ws.Cells.Item(12, 29).Value = 7  # AC12: 7.5 -> 7
$ws.Cells.Item(12, 35).Value = 9  # AI12: 8.5 -> 9
$ws.Cells.Item(12, 36).Value = 9.5  # AJ12: 9 -> 9.5
$ws.Cells.Item(12, 37).Value = 19  # AK12: 17 -> 19
$ws.Cells.Item(12, 41).Value = 21  # AO12: 23 -> 21
$ws.Cells.Item(12, 51).Value = 12  # AY12: 11 -> 12
$ws.Cells.Item(12, 52).Value = 26  # AZ12: 23 -> 26

# Row 14
$ws.Cells.Item(14, 7).Value = 2.55  # G14: 2.5 -> 2.55
$ws.Cells.Item(14, 9).Value = 2.55  # I14: 2.6 -> 2.55
$ws.Cells.Item(14, 10).Value = 3.2  # J14: 3.1 -> 3.2
$ws.Cells.Item(14, 17).Value = 1.73  # Q14: 1.8 -> 1.73
$ws.Cells.Item(14, 18).Value = 2.08  # R14: 2 -> 2.08
$ws.Cells.Item(14, 19).Value = 1.33  # S14: 1.36 -> 1.33
$ws.Cells.Item(14, 20).Value = 3.25  # T14: 3 -> 3.25
$ws.Cells.Item(14, 26).Value = 26  # Z14: 23 -> 26
$ws.Cells.Item(14, 29).Value = 13  # AC14: 12 -> 13
$ws.Cells.Item(14, 38).Value = 19  # AL14: 21 -> 19
$ws.Cells.Item(14, 46).Value = 3.25  # AT14: 3 -> 3.25
$ws.Cells.Item(14, 51).Value = 13  # AY14: 15 -> 13

# Row 15
$ws.Cells.Item(15, 14).Value = 9  # N15: 8.5 -> 9

# Row 16
$ws.Cells.Item(16, 7).Value = 1.36  # G16: 1.33 -> 1.36
$ws.Cells.Item(16, 9).Value = 8.5  # I16: 9 -> 8.5
$ws.Cells.Item(16, 17).Value = 1.67  # Q16: 1.7 -> 1.67
$ws.Cells.Item(16, 18).Value = 2.15  # R16: 2.1 -> 2.15
$ws.Cells.Item(16, 19).Value = 1.3  # S16: 1.33 -> 1.3
$ws.Cells.Item(16, 20).Value = 3.4  # T16: 3.25 -> 3.4
$ws.Cells.Item(16, 23).Value = 7.5  # W16: 7 -> 7.5
$ws.Cells.Item(16, 25).Value = 9  # Y16: 8.5 -> 9
$ws.Cells.Item(16, 27).Value = 11  # AA16: 12 -> 11
$ws.Cells.Item(16, 34).Value = 19  # AH16: 21 -> 19
$ws.Cells.Item(16, 36).Value = 23  # AJ16: 26 -> 23
$ws.Cells.Item(16, 37).Value = 81  # AK16: 101 -> 81
$ws.Cells.Item(16, 40).Value = 3.4  # AN16: 3.25 -> 3.4
$ws.Cells.Item(16, 42).Value = 17  # AP16: 19 -> 17
$ws.Cells.Item(16, 46).Value = 3.4  # AT16: 3.25 -> 3.4
$ws.Cells.Item(16, 47).Value = 9  # AU16: 9.5 -> 9
$ws.Cells.Item(16, 55).Value = 301  # BC16: 351 -> 301

# Row 17
$ws.Cells.Item(17, 15).Value = 1.18  # O17: 1.2 -> 1.18
$ws.Cells.Item(17, 16).Value = 4.5  # P17: 4.33 -> 4.5
$ws.Cells.Item(17, 17).Value = 1.65  # Q17: 1.67 -> 1.65
$ws.Cells.Item(17, 18).Value = 2.2  # R17: 2.15 -> 2.2

# Row 18
$ws.Cells.Item(18, 8).Value = 3.3  # H18: 3.25 -> 3.3
$ws.Cells.Item(18, 11).Value = 2.25  # K18: 2.2 -> 2.25
$ws.Cells.Item(18, 12).Value = 3.2  # L18: 3.25 -> 3.2
$ws.Cells.Item(18, 14).Value = 13  # N18: 12 -> 13
$ws.Cells.Item(18, 15).Value = 1.22  # O18: 1.25 -> 1.22
$ws.Cells.Item(18, 16).Value = 4  # P18: 3.75 -> 4
$ws.Cells.Item(18, 17).Value = 1.8  # Q18: 1.83 -> 1.8
$ws.Cells.Item(18, 18).Value = 2  # R18: 2.03 -> 2
$ws.Cells.Item(18, 45).Value = 126  # AS18: 151 -> 126
$ws.Cells.Item(18, 54).Value = 51  # BB18: 67 -> 51
$ws.Cells.Item(18, 55).Value = 126  # BC18: 151 -> 126

# Row 19
$ws.Cells.Item(19, 15).Value = 1.5  # O19: 1.44 -> 1.5
$ws.Cells.Item(19, 16).Value = 2.5  # P19: 2.63 -> 2.5
$ws.Cells.Item(19, 19).Value = 1.57  # S19: 1.54 -> 1.57

# Row 37
$ws.Cells.Item(37, 15).Value = 1.25  # O37: 1.29 -> 1.25
$ws.Cells.Item(37, 16).Value = 3.75  # P37: 3.5 -> 3.75
$ws.Cells.Item(37, 17).Value = 1.93  # Q37: 1.95 -> 1.93
$ws.Cells.Item(37, 18).Value = 1.93  # R37: 1.9 -> 1.93
